# Weekly update: a new daily price record is inserted as row 47 of the
# "Poroto granado" sheet, pushing the existing rows 47-107 down to 48-108.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 47 (shifts 47:107 -> 48:108)
$ws.Rows("47").Insert()

# Populate the new row 47 with the new weekly record
$ws.Range("A47").Value = 5
$ws.Range("B47").Value = "Macroferia Regional de Talca"
$ws.Range("C47").Value = "Maule"
$ws.Range("D47").Value = 44601
$ws.Range("E47").Value = 7
$ws.Range("F47").Value = 100112030
$ws.Range("G47").Value = "Poroto granado"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 300
$ws.Range("K47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("M47").Value = 20000
$ws.Range("N47").Value = "$/saco 25 kilos"
$ws.Range("O47").Value = "Región del Maule"
$ws.Range("P47").Value = 800
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"

Write-Host "Row 47 inserted and populated"
